$d = $word.ActiveDocument

# 1. Fix "whole page" leading space
$d.Content.Find.Execute(" whole page", $true, $false, $false, $false, $false, $true, 1, $false, "whole page", 2)

# 2. Fix "WE ARE USIN the " -> "WE ARE USING the "
$d.Content.Find.Execute("WE ARE USIN the", $true, $false, $false, $false, $false, $true, 1, $false, "WE ARE USING the", 2)
